# Update natmi TPM output values (Col2a1-Ddr1), rows 2-17,
# columns E:J and M:T, per new TPM recomputation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns E:J (rows 2-17)
$arrEJ = New-Object 'object[,]' 16,6
$arrEJ[0,0] = [double]"3"
$arrEJ[0,1] = [double]"1"
$arrEJ[0,2] = [double]"0.06437833333333333"
$arrEJ[0,3] = [double]"0.193135"
$arrEJ[0,4] = [double]"0.109187438766332"
$arrEJ[0,5] = [double]"0.109187438766332"
$arrEJ[1,0] = [double]"3"
$arrEJ[1,1] = [double]"1"
$arrEJ[1,2] = [double]"0.06437833333333333"
$arrEJ[1,3] = [double]"0.193135"
$arrEJ[1,4] = [double]"0.109187438766332"
$arrEJ[1,5] = [double]"0.109187438766332"
$arrEJ[2,0] = [double]"3"
$arrEJ[2,1] = [double]"1"
$arrEJ[2,2] = [double]"0.06437833333333333"
$arrEJ[2,3] = [double]"0.193135"
$arrEJ[2,4] = [double]"0.109187438766332"
$arrEJ[2,5] = [double]"0.109187438766332"
$arrEJ[3,0] = [double]"3"
$arrEJ[3,1] = [double]"1"
$arrEJ[3,2] = [double]"0.06437833333333333"
$arrEJ[3,3] = [double]"0.193135"
$arrEJ[3,4] = [double]"0.109187438766332"
$arrEJ[3,5] = [double]"0.109187438766332"
$arrEJ[4,0] = [double]"3"
$arrEJ[4,1] = [double]"1"
$arrEJ[4,2] = [double]"0.4788196666666666"
$arrEJ[4,3] = [double]"1.436459"
$arrEJ[4,4] = [double]"0.8120914339857952"
$arrEJ[4,5] = [double]"0.8120914339857951"
$arrEJ[5,0] = [double]"3"
$arrEJ[5,1] = [double]"1"
$arrEJ[5,2] = [double]"0.4788196666666666"
$arrEJ[5,3] = [double]"1.436459"
$arrEJ[5,4] = [double]"0.8120914339857952"
$arrEJ[5,5] = [double]"0.8120914339857951"
$arrEJ[6,0] = [double]"3"
$arrEJ[6,1] = [double]"1"
$arrEJ[6,2] = [double]"0.4788196666666666"
$arrEJ[6,3] = [double]"1.436459"
$arrEJ[6,4] = [double]"0.8120914339857952"
$arrEJ[6,5] = [double]"0.8120914339857951"
$arrEJ[7,0] = [double]"3"
$arrEJ[7,1] = [double]"1"
$arrEJ[7,2] = [double]"0.4788196666666666"
$arrEJ[7,3] = [double]"1.436459"
$arrEJ[7,4] = [double]"0.8120914339857952"
$arrEJ[7,5] = [double]"0.8120914339857951"
$arrEJ[8,0] = [double]"3"
$arrEJ[8,1] = [double]"1"
$arrEJ[8,2] = [double]"0.042481"
$arrEJ[8,3] = [double]"0.127443"
$arrEJ[8,4] = [double]"0.07204895414449818"
$arrEJ[8,5] = [double]"0.07204895414449818"
$arrEJ[9,0] = [double]"3"
$arrEJ[9,1] = [double]"1"
$arrEJ[9,2] = [double]"0.042481"
$arrEJ[9,3] = [double]"0.127443"
$arrEJ[9,4] = [double]"0.07204895414449818"
$arrEJ[9,5] = [double]"0.07204895414449818"
$arrEJ[10,0] = [double]"3"
$arrEJ[10,1] = [double]"1"
$arrEJ[10,2] = [double]"0.042481"
$arrEJ[10,3] = [double]"0.127443"
$arrEJ[10,4] = [double]"0.07204895414449818"
$arrEJ[10,5] = [double]"0.07204895414449818"
$arrEJ[11,0] = [double]"3"
$arrEJ[11,1] = [double]"1"
$arrEJ[11,2] = [double]"0.042481"
$arrEJ[11,3] = [double]"0.127443"
$arrEJ[11,4] = [double]"0.07204895414449818"
$arrEJ[11,5] = [double]"0.07204895414449818"
$arrEJ[12,0] = [double]"1"
$arrEJ[12,1] = [double]"0.3333333333333333"
$arrEJ[12,2] = [double]"0.003934"
$arrEJ[12,3] = [double]"0.011802"
$arrEJ[12,4] = [double]"0.006672173103374587"
$arrEJ[12,5] = [double]"0.006672173103374586"
$arrEJ[13,0] = [double]"1"
$arrEJ[13,1] = [double]"0.3333333333333333"
$arrEJ[13,2] = [double]"0.003934"
$arrEJ[13,3] = [double]"0.011802"
$arrEJ[13,4] = [double]"0.006672173103374587"
$arrEJ[13,5] = [double]"0.006672173103374586"
$arrEJ[14,0] = [double]"1"
$arrEJ[14,1] = [double]"0.3333333333333333"
$arrEJ[14,2] = [double]"0.003934"
$arrEJ[14,3] = [double]"0.011802"
$arrEJ[14,4] = [double]"0.006672173103374587"
$arrEJ[14,5] = [double]"0.006672173103374586"
$arrEJ[15,0] = [double]"1"
$arrEJ[15,1] = [double]"0.3333333333333333"
$arrEJ[15,2] = [double]"0.003934"
$arrEJ[15,3] = [double]"0.011802"
$arrEJ[15,4] = [double]"0.006672173103374587"
$arrEJ[15,5] = [double]"0.006672173103374586"
$ws.Range("E2:J17").Value = $arrEJ

# Columns M:T (rows 2-17)
$arrMT = New-Object 'object[,]' 16,8
$arrMT[0,0] = [double]"0.3360566666666667"
$arrMT[0,1] = [double]"1.00817"
$arrMT[0,2] = [double]"0.01570866217798777"
$arrMT[0,3] = [double]"0.01570866217798777"
$arrMT[0,4] = [double]"0.02163476810555556"
$arrMT[0,5] = [double]"0.19471291295"
$arrMT[0,6] = [double]"0.001715188589660036"
$arrMT[0,7] = [double]"0.001715188589660036"
$arrMT[1,0] = [double]"2.338622"
$arrMT[1,1] = [double]"7.015866"
$arrMT[1,2] = [double]"0.109316751024163"
$arrMT[1,3] = [double]"0.1093167510241629"
$arrMT[1,4] = [double]"0.1505565866566667"
$arrMT[1,5] = [double]"1.35500927991"
$arrMT[1,6] = [double]"0.01193601605858516"
$arrMT[1,7] = [double]"0.01193601605858516"
$arrMT[2,0] = [double]"18.491866"
$arrMT[2,1] = [double]"55.47559800000001"
$arrMT[2,2] = [double]"0.864385399390831"
$arrMT[2,3] = [double]"0.864385399390831"
$arrMT[2,4] = [double]"1.190475513303333"
$arrMT[2,5] = [double]"10.71427961973"
$arrMT[2,6] = [double]"0.09438002786649782"
$arrMT[2,7] = [double]"0.09438002786649782"
$arrMT[3,0] = [double]"0.2265353333333333"
$arrMT[3,1] = [double]"0.6796059999999999"
$arrMT[3,2] = [double]"0.01058918740701822"
$arrMT[3,3] = [double]"0.01058918740701822"
$arrMT[3,4] = [double]"0.01458396720111111"
$arrMT[3,5] = [double]"0.13125570481"
$arrMT[3,6] = [double]"0.001156206251589016"
$arrMT[3,7] = [double]"0.001156206251589016"
$arrMT[4,0] = [double]"0.3360566666666667"
$arrMT[4,1] = [double]"1.00817"
$arrMT[4,2] = [double]"0.01570866217798777"
$arrMT[4,3] = [double]"0.01570866217798777"
$arrMT[4,4] = [double]"0.1609105411144444"
$arrMT[4,5] = [double]"1.44819487003"
$arrMT[4,6] = [double]"0.01275686999412051"
$arrMT[4,7] = [double]"0.01275686999412051"
$arrMT[5,0] = [double]"2.338622"
$arrMT[5,1] = [double]"7.015866"
$arrMT[5,2] = [double]"0.109316751024163"
$arrMT[5,3] = [double]"0.1093167510241629"
$arrMT[5,4] = [double]"1.119778206499333"
$arrMT[5,5] = [double]"10.078003858494"
$arrMT[5,6] = [double]"0.08877519709788063"
$arrMT[5,7] = [double]"0.08877519709788062"
$arrMT[6,0] = [double]"18.491866"
$arrMT[6,1] = [double]"55.47559800000001"
$arrMT[6,2] = [double]"0.864385399390831"
$arrMT[6,3] = [double]"0.864385399390831"
$arrMT[6,4] = [double]"8.854269114164667"
$arrMT[6,5] = [double]"79.68842202748201"
$arrMT[6,6] = [double]"0.7019599785076843"
$arrMT[6,7] = [double]"0.7019599785076842"
$arrMT[7,0] = [double]"0.2265353333333333"
$arrMT[7,1] = [double]"0.6796059999999999"
$arrMT[7,2] = [double]"0.01058918740701822"
$arrMT[7,3] = [double]"0.01058918740701822"
$arrMT[7,4] = [double]"0.1084695727948889"
$arrMT[7,5] = [double]"0.9762261551539999"
$arrMT[7,6] = [double]"0.00859938838610975"
$arrMT[7,7] = [double]"0.008599388386109748"
$arrMT[8,0] = [double]"0.3360566666666667"
$arrMT[8,1] = [double]"1.00817"
$arrMT[8,2] = [double]"0.01570866217798777"
$arrMT[8,3] = [double]"0.01570866217798777"
$arrMT[8,4] = [double]"0.01427602325666667"
$arrMT[8,5] = [double]"0.12848420931"
$arrMT[8,6] = [double]"0.001131792680933254"
$arrMT[8,7] = [double]"0.001131792680933254"
$arrMT[9,0] = [double]"2.338622"
$arrMT[9,1] = [double]"7.015866"
$arrMT[9,2] = [double]"0.109316751024163"
$arrMT[9,3] = [double]"0.1093167510241629"
$arrMT[9,4] = [double]"0.09934700118199999"
$arrMT[9,5] = [double]"0.894123010638"
$arrMT[9,6] = [double]"0.007876157581765441"
$arrMT[9,7] = [double]"0.00787615758176544"
$arrMT[10,0] = [double]"18.491866"
$arrMT[10,1] = [double]"55.47559800000001"
$arrMT[10,2] = [double]"0.864385399390831"
$arrMT[10,3] = [double]"0.864385399390831"
$arrMT[10,4] = [double]"0.785552959546"
$arrMT[10,5] = [double]"7.069976635914001"
$arrMT[10,6] = [double]"0.06227806400388373"
$arrMT[10,7] = [double]"0.06227806400388373"
$arrMT[11,0] = [double]"0.2265353333333333"
$arrMT[11,1] = [double]"0.6796059999999999"
$arrMT[11,2] = [double]"0.01058918740701822"
$arrMT[11,3] = [double]"0.01058918740701822"
$arrMT[11,4] = [double]"0.009623447495333332"
$arrMT[11,5] = [double]"0.08661102745799999"
$arrMT[11,6] = [double]"0.0007629398779157531"
$arrMT[11,7] = [double]"0.0007629398779157531"
$arrMT[12,0] = [double]"0.3360566666666667"
$arrMT[12,1] = [double]"1.00817"
$arrMT[12,2] = [double]"0.01570866217798777"
$arrMT[12,3] = [double]"0.01570866217798777"
$arrMT[12,4] = [double]"0.001322046926666667"
$arrMT[12,5] = [double]"0.01189842234"
$arrMT[12,6] = [double]"0.0001048109132739677"
$arrMT[12,7] = [double]"0.0001048109132739676"
$arrMT[13,0] = [double]"2.338622"
$arrMT[13,1] = [double]"7.015866"
$arrMT[13,2] = [double]"0.109316751024163"
$arrMT[13,3] = [double]"0.1093167510241629"
$arrMT[13,4] = [double]"0.009200138948"
$arrMT[13,5] = [double]"0.08280125053199999"
$arrMT[13,6] = [double]"0.0007293802859317163"
$arrMT[13,7] = [double]"0.0007293802859317162"
$arrMT[14,0] = [double]"18.491866"
$arrMT[14,1] = [double]"55.47559800000001"
$arrMT[14,2] = [double]"0.864385399390831"
$arrMT[14,3] = [double]"0.864385399390831"
$arrMT[14,4] = [double]"0.072747000844"
$arrMT[14,5] = [double]"0.6547230075960001"
$arrMT[14,6] = [double]"0.005767329012765203"
$arrMT[14,7] = [double]"0.005767329012765202"
$arrMT[15,0] = [double]"0.2265353333333333"
$arrMT[15,1] = [double]"0.6796059999999999"
$arrMT[15,2] = [double]"0.01058918740701822"
$arrMT[15,3] = [double]"0.01058918740701822"
$arrMT[15,4] = [double]"0.0008911900013333332"
$arrMT[15,5] = [double]"0.008020710011999999"
$arrMT[15,6] = [double]"7.065289140369984E-05"
$arrMT[15,7] = [double]"7.065289140369982E-05"
$ws.Range("M2:T17").Value = $arrMT
